# "aggiornamento fino a 13/03" - append 4 more daily rows (2021-05-10 .. 2021-05-13,
# serials 44326-44329) below the existing data, mirroring the formatting of the
# last existing row (A251) and zero-filling the B/C/D metric columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 251
$newDates = @(44326, 44327, 44328, 44329)

# Copy the date cell's formatting (style) from the last existing row down onto
# the new rows, then overwrite values explicitly.
$ws.Range("A$lastRow").Copy() | Out-Null
$ws.Range("A" + ($lastRow + 1) + ":A" + ($lastRow + $newDates.Count)).PasteSpecial(-4122) | Out-Null

$r = $lastRow + 1
foreach ($d in $newDates) {
    $ws.Cells.Item($r, 1).Value = $d
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
    $r++
}
